# Translate the "typology" glossary (column I) of the "all law" sheet from
# English into Italian, and add a small glossary lookup sheet ("Sheet1")
# that maps the English terms to their Italian translations.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("all law")

# Translate the typology values (column I) for each law record.
# Row 1 is the header ("range") and is left untouched.
$ws.Range("I2:I339").Value   = "cantonale"
$ws.Range("I340:I347").Value = "federale"
$ws.Range("I348:I356").Value = "intercantonale"
$ws.Range("I357:I358").Value = "internazionale"

# Turn on the auto-filter for the full table.
$ws.Range("A1:I358").AutoFilter()

# Update the (hidden) filter-database defined name so it covers the whole
# table, including the newly added "typology" column / rows.
$filterName = $wb.Names.Item(1)
$filterName.RefersTo = "='all law'!`$A`$1:`$I`$358"

# Add a new worksheet right after "all law" holding the EN -> IT glossary.
$glossary = $wb.Worksheets.Add($null, $ws)
$glossary.Name = "Sheet1"

$glossary.Range("A1").Value = "cantonal"
$glossary.Range("B1").Value = "cantonale"
$glossary.Range("A2").Value = "federal"
$glossary.Range("B2").Value = "federale"
$glossary.Range("A3").Value = "intercantonal"
$glossary.Range("B3").Value = "intercantonale"
$glossary.Range("A4").Value = "international"
$glossary.Range("B4").Value = "internazionale"
$glossary.Range("B4").Select()

# Leave the focus back on the "all law" sheet, with column J selected (as
# was left by the author after finishing the translation pass).
$ws.Activate()
$ws.Columns.Item(10).Select()
